$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44483
$ws.Range("H2").Value = 'Camote'
$ws.Range("I2").Value = '1a nueva(o)'
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 550
$ws.Range("L2").Value = 580
$ws.Range("M2").Value = 565
$ws.Range("O2").Value = 'Perú'
$ws.Range("P2").Value = 565

$ws.Range("D3").Value = 44547
$ws.Range("H3").Value = 'Camote'
$ws.Range("I3").Value = '1a (cosecha)'
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 650
$ws.Range("M3").Value = 625
$ws.Range("O3").Value = 'Perú'
$ws.Range("P3").Value = 625

$ws.Range("D4").Value = 44547
$ws.Range("H4").Value = 'Camote'
$ws.Range("I4").Value = '2a nueva(o)'
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 550
$ws.Range("M4").Value = 525
$ws.Range("O4").Value = 'Perú'
$ws.Range("P4").Value = 525

$ws.Range("D5").Value = 44211
$ws.Range("H5").Value = 'Camote'
$ws.Range("I5").Value = '1a nueva(o)'
$ws.Range("J5").Value = 1600
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 550
$ws.Range("M5").Value = 525
$ws.Range("O5").Value = 'Región de O''Higgins'
$ws.Range("P5").Value = 525

$ws.Range("D6").Value = 44575
$ws.Range("H6").Value = 'Camote'
$ws.Range("I6").Value = '1a nueva(o)'
$ws.Range("J6").Value = 1300
$ws.Range("K6").Value = 500
$ws.Range("L6").Value = 550
$ws.Range("M6").Value = 525
$ws.Range("O6").Value = 'Región de O''Higgins'
$ws.Range("P6").Value = 525

$ws.Range("D7").Value = 44630
$ws.Range("H7").Value = 'Camote'
$ws.Range("I7").Value = '1a (cosecha)'
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 450
$ws.Range("L7").Value = 480
$ws.Range("M7").Value = 465
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 465

$ws.Range("D8").Value = 44301
$ws.Range("H8").Value = 'Camote'
$ws.Range("I8").Value = '2a nueva(o)'
$ws.Range("J8").Value = 1200
$ws.Range("K8").Value = 400
$ws.Range("L8").Value = 430
$ws.Range("M8").Value = 415
$ws.Range("O8").Value = 'Provincia de Melipilla'
$ws.Range("P8").Value = 415

$ws.Range("D9").Value = 44530
$ws.Range("H9").Value = 'Camote'
$ws.Range("I9").Value = '2a nueva(o)'
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 480
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = 490
$ws.Range("O9").Value = 'Perú'
$ws.Range("P9").Value = 490

$ws.Range("D10").Value = 44620
$ws.Range("H10").Value = 'Camote'
$ws.Range("I10").Value = '1a (cosecha)'
$ws.Range("J10").Value = 1200
$ws.Range("K10").Value = 480
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = 490
$ws.Range("O10").Value = 'Región de O''Higgins'
$ws.Range("P10").Value = 490

$ws.Range("D11").Value = 44175
$ws.Range("H11").Value = 'Camote'
$ws.Range("I11").Value = '1a nueva(o)'
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1550
$ws.Range("M11").Value = 1525
$ws.Range("O11").Value = 'Perú'
$ws.Range("P11").Value = 1525

$ws.Range("D12").Value = 44399
$ws.Range("H12").Value = 'Camote'
$ws.Range("I12").Value = '1a (guarda)'
$ws.Range("J12").Value = 800
$ws.Range("K12").Value = 450
$ws.Range("L12").Value = 480
$ws.Range("M12").Value = 465
$ws.Range("O12").Value = 'Provincia de Melipilla'
$ws.Range("P12").Value = 465

$ws.Range("D13").Value = 44691
$ws.Range("H13").Value = 'Camote'
$ws.Range("I13").Value = '1a (cosecha)'
$ws.Range("J13").Value = 700
$ws.Range("K13").Value = 580
$ws.Range("L13").Value = 600
$ws.Range("M13").Value = 590
$ws.Range("O13").Value = 'Región de O''Higgins'
$ws.Range("P13").Value = 590

$ws.Range("D14").Value = 44736
$ws.Range("H14").Value = 'Camote'
$ws.Range("I14").Value = '1a (guarda)'
$ws.Range("J14").Value = 900
$ws.Range("K14").Value = 780
$ws.Range("L14").Value = 800
$ws.Range("M14").Value = 790
$ws.Range("O14").Value = 'Región de O''Higgins'
$ws.Range("P14").Value = 790

$ws.Range("D15").Value = 44670
$ws.Range("H15").Value = 'Camote'
$ws.Range("I15").Value = '1a (cosecha)'
$ws.Range("J15").Value = 1200
$ws.Range("K15").Value = 400
$ws.Range("L15").Value = 430
$ws.Range("M15").Value = 415
$ws.Range("O15").Value = 'Región de O''Higgins'
$ws.Range("P15").Value = 415

$ws.Range("D16").Value = 44453
$ws.Range("H16").Value = 'Camote'
$ws.Range("I16").Value = '1a nueva(o)'
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 630
$ws.Range("L16").Value = 650
$ws.Range("M16").Value = 640
$ws.Range("O16").Value = 'Perú'
$ws.Range("P16").Value = 640

$ws.Range("D17").Value = 44601
$ws.Range("H17").Value = 'Camote'
$ws.Range("I17").Value = '2a (cosecha)'
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 400
$ws.Range("L17").Value = 450
$ws.Range("M17").Value = 425
$ws.Range("O17").Value = 'Región de O''Higgins'
$ws.Range("P17").Value = 425

$ws.Range("D18").Value = 44179
$ws.Range("H18").Value = 'Camote'
$ws.Range("I18").Value = '1a nueva(o)'
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1350
$ws.Range("L18").Value = 1400
$ws.Range("M18").Value = 1375
$ws.Range("O18").Value = 'Perú'
$ws.Range("P18").Value = 1375

$ws.Range("D19").Value = 44428
$ws.Range("H19").Value = 'Camote'
$ws.Range("I19").Value = '1a nueva(o)'
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 580
$ws.Range("L19").Value = 600
$ws.Range("M19").Value = 590
$ws.Range("O19").Value = 'Perú'
$ws.Range("P19").Value = 590

$ws.Range("D20").Value = 44204
$ws.Range("H20").Value = 'Camote'
$ws.Range("I20").Value = '2a nueva(o)'
$ws.Range("J20").Value = 1600
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 550
$ws.Range("M20").Value = 525
$ws.Range("O20").Value = 'Región del Maule'
$ws.Range("P20").Value = 525

$ws.Range("D21").Value = 44469
$ws.Range("H21").Value = 'Camote'
$ws.Range("I21").Value = '1a nueva(o)'
$ws.Range("J21").Value = 1200
$ws.Range("K21").Value = 600
$ws.Range("L21").Value = 650
$ws.Range("M21").Value = 625
$ws.Range("O21").Value = 'Perú'
$ws.Range("P21").Value = 625

$ws.Range("D22").Value = 44349
$ws.Range("H22").Value = 'Pachia'
$ws.Range("I22").Value = '1a nueva(o)'
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 730
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = 740
$ws.Range("O22").Value = 'Perú'
$ws.Range("P22").Value = 740

$ws.Range("D23").Value = 44650
$ws.Range("H23").Value = 'Camote'
$ws.Range("I23").Value = '2a (cosecha)'
$ws.Range("J23").Value = 1300
$ws.Range("K23").Value = 400
$ws.Range("L23").Value = 430
$ws.Range("M23").Value = 415
$ws.Range("O23").Value = 'Región de O''Higgins'
$ws.Range("P23").Value = 415

$ws.Range("D24").Value = 44243
$ws.Range("H24").Value = 'Camote'
$ws.Range("I24").Value = '2a nueva(o)'
$ws.Range("J24").Value = 1600
$ws.Range("K24").Value = 450
$ws.Range("L24").Value = 480
$ws.Range("M24").Value = 465
$ws.Range("O24").Value = 'Región del Maule'
$ws.Range("P24").Value = 465

$ws.Range("D25").Value = 44201
$ws.Range("H25").Value = 'Camote'
$ws.Range("I25").Value = '1a nueva(o)'
$ws.Range("J25").Value = 1360
$ws.Range("K25").Value = 730
$ws.Range("L25").Value = 750
$ws.Range("M25").Value = 740
$ws.Range("O25").Value = 'Perú'
$ws.Range("P25").Value = 740

$ws.Range("D26").Value = 44231
$ws.Range("H26").Value = 'Camote'
$ws.Range("I26").Value = '1a nueva(o)'
$ws.Range("J26").Value = 1300
$ws.Range("K26").Value = 450
$ws.Range("L26").Value = 480
$ws.Range("M26").Value = 465
$ws.Range("O26").Value = 'Perú'
$ws.Range("P26").Value = 465

$ws.Range("D27").Value = 44490
$ws.Range("H27").Value = 'Camote'
$ws.Range("I27").Value = '1a nueva(o)'
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 480
$ws.Range("M27").Value = 465
$ws.Range("O27").Value = 'Perú'
$ws.Range("P27").Value = 465

$ws.Range("D28").Value = 44476
$ws.Range("H28").Value = 'Camote'
$ws.Range("I28").Value = '1a nueva(o)'
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 480
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 490
$ws.Range("O28").Value = 'Perú'
$ws.Range("P28").Value = 490

$ws.Range("D29").Value = 44685
$ws.Range("H29").Value = 'Camote'
$ws.Range("I29").Value = '1a (cosecha)'
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 680
$ws.Range("L29").Value = 700
$ws.Range("M29").Value = 690
$ws.Range("O29").Value = 'Región de O''Higgins'
$ws.Range("P29").Value = 690

$ws.Range("D30").Value = 44238
$ws.Range("H30").Value = 'Camote'
$ws.Range("I30").Value = '1a nueva(o)'
$ws.Range("J30").Value = 1250
$ws.Range("K30").Value = 430
$ws.Range("L30").Value = 450
$ws.Range("M30").Value = 440
$ws.Range("O30").Value = 'Perú'
$ws.Range("P30").Value = 440

$ws.Range("D31").Value = 44322
$ws.Range("H31").Value = 'Camote'
$ws.Range("I31").Value = '1a (cosecha)'
$ws.Range("J31").Value = 1200
$ws.Range("K31").Value = 350
$ws.Range("L31").Value = 400
$ws.Range("M31").Value = 375
$ws.Range("O31").Value = 'Región del Maule'
$ws.Range("P31").Value = 375

$ws.Range("D32").Value = 44665
$ws.Range("H32").Value = 'Camote'
$ws.Range("I32").Value = '1a (cosecha)'
$ws.Range("J32").Value = 1200
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 420
$ws.Range("M32").Value = 410
$ws.Range("O32").Value = 'Región de O''Higgins'
$ws.Range("P32").Value = 410
